# Creating and Reading Excel using Apache POI
# Add a new "Department Data" worksheet after the existing "Employee Data" sheet
# and populate it with department records.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Department Data"

$data = @(
    @("DeptNum", "DeptName", "DeptLocation"),
    @("'10", "Dept10", "India"),
    @("'20", "Dept20", "UK"),
    @("'30", "Dept30", "USA"),
    @("'40", "Dept40", "Japan"),
    @("'50", "Dept50", "Russia")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# The numeric-looking department numbers were entered with a leading
# apostrophe so Excel stores them as text (matching the source data's
# shared-string cell type). Clear the resulting "quote prefix" formatting
# so the cells keep the default style, only the text value/type remains.
$ws.Range("A1:C6").ClearFormats()
